$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the lat/lng correction formulas: the per-second offset multipliers
# change from 0.0000001 to 0.000001 (latitude, J/K) and 0.0001 (longitude, L/M).
$ws.Range("B2").Formula = "=H2+(J2*0.000001)-(K2*0.000001)"
$ws.Range("C2").Formula = "=I2+(L2*0.0001)-(M2*0.0001)"

# B3:B7 / C3:C7 share the same formula pattern down the column.
$ws.Range("B3:B7").Formula = "=H3+(J3*0.000001)-(K3*0.000001)"
$ws.Range("C3:C7").Formula = "=I3+(L3*0.0001)-(M3*0.0001)"

# New longitude-seconds input for the first location.
$ws.Range("L2").Value = 45

# Move the active selection from D12 to D10.
$ws.Range("D10").Select()
